# Refresh the coinranking.com snapshot cells touched by the scheduled
# GitHub Actions run ("Updated cryptos list on Fri Oct 20 18:14:39 UTC 2023").
#
# Price (D) and change (E) text is re-written for every row whose quote
# moved; two rows (44/45 and 50/51) also had their rank swap, so the coin
# name (B) and coinranking link (C) move along with the new D/E values.
#
# Price cells that look like a bare decimal (e.g. "0.999", "4.00",
# "0.0490") must stay plain text -- otherwise Excel's COM layer silently
# coerces them to doubles and the significant trailing zeros are lost.
# Flip the cell to text, write the string, then restore the "Normal"
# style so no stray number-format style sticks to the cell afterwards.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.471.15'
$ws.Range('E2').Value = '  +3.16%  '
$ws.Range('D3').Value = '1.605.58'
$ws.Range('E3').Value = '  +2.77%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.75'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.07%  '
$ws.Range('E6').Value = '  +6.74%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '26.83'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +7.66%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '43.53'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.251'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.75%  '
$ws.Range('E11').Value = '  +2.46%  '
$ws.Range('D13').Value = '1.836.11'
$ws.Range('E13').Value = '  +2.86%  '
$ws.Range('D14').Value = '1.627.87'
$ws.Range('E14').Value = '  +4.58%  '
$ws.Range('D15').Value = '29.492.61'
$ws.Range('E15').Value = '  +3.24%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.535'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.19%  '
$ws.Range('E17').Value = '  +2.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.10'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '241.67'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.64'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.03%  '
$ws.Range('D21').Value = '0.0₃0690'
$ws.Range('E21').Value = '  +1.81%  '
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.33%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.21'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.66%  '
$ws.Range('E25').Value = '  +0.41%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.41'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.57%  '
$ws.Range('E27').Value = '  +5.13%  '
$ws.Range('E28').Value = '  +3.55%  '
$ws.Range('E29').Value = '  +2.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0473'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.80%  '
$ws.Range('E32').Value = '  +0.93%  '
$ws.Range('E33').Value = '  +1.85%  '
$ws.Range('E34').Value = '  +4.65%  '
$ws.Range('D35').Value = '1.414.74'
$ws.Range('E35').Value = '  +2.01%  '
$ws.Range('E36').Value = '  +0.48%  '
$ws.Range('E37').Value = '  +3.52%  '
$ws.Range('E38').Value = '  +5.07%  '
$ws.Range('E39').Value = '  +0.16%  '
$ws.Range('E40').Value = '  +2.82%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.537'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.63%  '
$ws.Range('E42').Value = '  +1.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0490'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.82%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.798'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.41%  '
$ws.Range('B45').Value = 'BitcoinSV'
$ws.Range('C45').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '52.90'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +22.39%  '
$ws.Range('E46').Value = '  -0.12%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '65.70'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.93%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.28'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.22%  '
$ws.Range('D49').Value = '1.746.59'
$ws.Range('E49').Value = '  +3.01%  '
$ws.Range('B50').Value = 'WEMIXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.852'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.04%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '86.89'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.13%  '
